$p = $ppt.ActivePresentation

# Each of the 6 existing slides currently holds a single free-floating
# TextBox ('TextBox 1'). We replace it with a proper Title + Content
# placeholder pair (ppLayoutText / CustomLayout 'Title and Content'),
# matching the target deck's shape structure, then fill in the new
# generated LinkedIn-carousel copy. We insert the replacement slide
# right before the original, fill it in, then delete the original so
# the overall slide order (1..6) is preserved.

# --- Slide 1 ---
$ns = $p.Slides.Add(1, 2)
$titleTf = $ns.Shapes.Item(1).TextFrame
$titleTf.TextRange.Font.Size = 28
$titleTf.TextRange.Font.Bold = $true
$titleTf.TextRange.Text = 'Elcogen Enters India and APAC with Green Hydrogen Technology - Saur Energy'

$bodyTf = $ns.Shapes.Item(2).TextFrame
$bodyLines = @('⚠️ Image missing', 'Elcogen has expanded its operations into India and the Asia-Pacific (APAC) region.', 'The company focuses on green hydrogen technology.')
$bodyTf.TextRange.Text = [string]::Join("`r", $bodyLines)
for ($li = 2; $li -le $bodyLines.Count; $li++) {
  $para = $bodyTf.TextRange.Paragraphs($li, 1)
  $para.IndentLevel = 2
  $para.Font.Size = 16
}

$p.Slides.Item(2).Delete()

# --- Slide 2 ---
$ns = $p.Slides.Add(2, 2)
$titleTf = $ns.Shapes.Item(1).TextFrame
$titleTf.TextRange.Font.Size = 28
$titleTf.TextRange.Font.Bold = $true
$titleTf.TextRange.Text = 'Elcogen Enters India and APAC with Green Hydrogen Technology - Saur Energy'

$bodyTf = $ns.Shapes.Item(2).TextFrame
$bodyLines = @('⚠️ Image missing', 'The expansion is part of Elcogen''s strategy to enhance its market presence in emerging economies.', 'Green hydrogen technology is a key component in the transition to sustainable energy.')
$bodyTf.TextRange.Text = [string]::Join("`r", $bodyLines)
for ($li = 2; $li -le $bodyLines.Count; $li++) {
  $para = $bodyTf.TextRange.Paragraphs($li, 1)
  $para.IndentLevel = 2
  $para.Font.Size = 16
}

$p.Slides.Item(3).Delete()

# --- Slide 3 ---
$ns = $p.Slides.Add(3, 2)
$titleTf = $ns.Shapes.Item(1).TextFrame
$titleTf.TextRange.Font.Size = 28
$titleTf.TextRange.Font.Bold = $true
$titleTf.TextRange.Text = 'Elcogen Enters India and APAC with Green Hydrogen Technology - Saur Energy'

$bodyTf = $ns.Shapes.Item(2).TextFrame
$bodyLines = @('⚠️ Image missing', 'Elcogen aims to leverage the growing demand for clean energy solutions in the APAC region.', 'The company is positioned to contribute to the reduction of carbon emissions.')
$bodyTf.TextRange.Text = [string]::Join("`r", $bodyLines)
for ($li = 2; $li -le $bodyLines.Count; $li++) {
  $para = $bodyTf.TextRange.Paragraphs($li, 1)
  $para.IndentLevel = 2
  $para.Font.Size = 16
}

$p.Slides.Item(4).Delete()

# --- Slide 4 ---
$ns = $p.Slides.Add(4, 2)
$titleTf = $ns.Shapes.Item(1).TextFrame
$titleTf.TextRange.Font.Size = 28
$titleTf.TextRange.Font.Bold = $true
$titleTf.TextRange.Text = 'Elcogen Enters India and APAC with Green Hydrogen Technology - Saur Energy'

$bodyTf = $ns.Shapes.Item(2).TextFrame
$bodyLines = @('⚠️ Image missing', 'The entry into India aligns with the country''s goals for renewable energy adoption.', 'Elcogen''s technology supports the development of hydrogen as a clean fuel source.')
$bodyTf.TextRange.Text = [string]::Join("`r", $bodyLines)
for ($li = 2; $li -le $bodyLines.Count; $li++) {
  $para = $bodyTf.TextRange.Paragraphs($li, 1)
  $para.IndentLevel = 2
  $para.Font.Size = 16
}

$p.Slides.Item(5).Delete()

# --- Slide 5 ---
$ns = $p.Slides.Add(5, 2)
$titleTf = $ns.Shapes.Item(1).TextFrame
$titleTf.TextRange.Font.Size = 28
$titleTf.TextRange.Font.Bold = $true
$titleTf.TextRange.Text = 'Elcogen Enters India and APAC with Green Hydrogen Technology - Saur Energy'

$bodyTf = $ns.Shapes.Item(2).TextFrame
$bodyLines = @('⚠️ Image missing', 'The initiative reflects a broader trend of international companies investing in India''s energy sector.', 'Elcogen''s technology is expected to play a role in various industrial applications.')
$bodyTf.TextRange.Text = [string]::Join("`r", $bodyLines)
for ($li = 2; $li -le $bodyLines.Count; $li++) {
  $para = $bodyTf.TextRange.Paragraphs($li, 1)
  $para.IndentLevel = 2
  $para.Font.Size = 16
}

$p.Slides.Item(6).Delete()

# --- Slide 6 ---
$ns = $p.Slides.Add(6, 2)
$titleTf = $ns.Shapes.Item(1).TextFrame
$titleTf.TextRange.Font.Size = 28
$titleTf.TextRange.Font.Bold = $true
$titleTf.TextRange.Text = 'Elcogen Enters India and APAC with Green Hydrogen Technology - Saur Energy'

$bodyTf = $ns.Shapes.Item(2).TextFrame
$bodyLines = @('⚠️ Image missing', 'The expansion is anticipated to create new opportunities for collaboration in the region.', 'Elcogen''s move is significant in the context of global efforts to achieve net-zero emissions.')
$bodyTf.TextRange.Text = [string]::Join("`r", $bodyLines)
for ($li = 2; $li -le $bodyLines.Count; $li++) {
  $para = $bodyTf.TextRange.Paragraphs($li, 1)
  $para.IndentLevel = 2
  $para.Font.Size = 16
}

$p.Slides.Item(7).Delete()

